$wb = $excel.ActiveWorkbook

# --- Work on "Iteration #2": adjust scroll position (cosmetic) ---
$ws2 = $wb.Worksheets.Item("Iteration #2")
$ws2.Activate()
$win2 = $excel.ActiveWindow
$win2.ScrollRow = 10
$win2.ScrollColumn = 1

# --- Work on "Iteration #3": the sheet that receives new iteration data ---
$ws3 = $wb.Worksheets.Item("Iteration #3")
$ws3.Activate()

# Row 14: first new entry
$ws3.Range("A14").Value2 = 42814
$ws3.Range("B14").Value = "Présentation des projets."
$ws3.Range("C14").Value = 3

# Row 15: second new entry
$ws3.Range("A15").Value2 = 42815
$ws3.Range("B15").Value = "Modifications de la structure du code terminées."
$ws3.Range("C15").Value = 2

# Make A15's date format match A14 (copy format only, avoids creating a new style)
$ws3.Range("A14").Copy()
$ws3.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Center the "Travail (h)" column for the data rows
$ws3.Range("C14:C36").HorizontalAlignment = -4108

# Scroll / selection state
$win3 = $excel.ActiveWindow
$win3.ScrollRow = 11
$win3.ScrollColumn = 1
$ws3.Range("B16:B17").Select()

$wb.Save()
